# Updates the cryptos sheet with refreshed price/volume figures.
#
# Price (column D) and Volume(1h) (column E) are stored as plain TEXT
# cells in the source workbook, not numbers - values like "27.739.53"
# (multiple dots) or "15.50" / "0.5200" (meaningful trailing zeros) and
# the padded "  +3.12%  " strings must come through byte-for-byte.
# Assigning a numeric-looking string straight to Range.Value lets Excel
# "helpfully" reinterpret it as a real number (dropping the trailing
# zero, etc.), so for any new Price value that could parse as a plain
# decimal number we force Text format on the cell first, assign the
# value, then clear the formatting again so the cell is left exactly as
# it started (General format, no explicit style) - only its text content
# has changed, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "27.739.53"
    "E2" = "  +3.12%  "
    "D3" = "1.865.12"
    "E3" = "  +2.89%  "
    "E4" = "  +3.11%  "
    "D5" = "324.65"
    "E5" = "  +3.87%  "
    "E6" = "  +2.76%  "
    "E7" = "  +2.79%  "
    "D8" = "0.3801"
    "E8" = "  +2.75%  "
    "D9" = "0.07463"
    "E9" = "  +2.80%  "
    "D10" = "0.8846"
    "E11" = "  +2.28%  "
    "D12" = "1.870.23"
    "E12" = "  -9.39%  "
    "D13" = "5.558"
    "E13" = "  +2.85%  "
    "D14" = "6.754"
    "E14" = "  +1.73%  "
    "D15" = "0.07221"
    "E15" = "  +4.12%  "
    "D16" = "83.95"
    "E16" = "  +3.83%  "
    "E17" = "  +2.67%  "
    "D18" = "0.000009108"
    "E18" = "  +3.08%  "
    "D19" = "1.034"
    "E19" = "  +2.76%  "
    "D20" = "15.50"
    "E20" = "  +1.99%  "
    "D21" = "27.762.13"
    "E21" = "  +3.06%  "
    "E22" = "  +2.17%  "
    "E23" = "  +4.20%  "
    "D24" = "1.958"
    "E24" = "  +3.96%  "
    "D25" = "158.19"
    "E25" = "  +2.54%  "
    "D26" = "18.88"
    "E26" = "  +2.79%  "
    "D27" = "1.999"
    "E27" = "  +4.02%  "
    "D28" = "5.308"
    "E28" = "  +1.20%  "
    "D29" = "117.60"
    "E29" = "  +2.55%  "
    "D30" = "0.09096"
    "E30" = "  +1.74%  "
    "D31" = "1.215"
    "E31" = "  +4.53%  "
    "D32" = "0.7710"
    "E32" = "  +3.63%  "
    "D33" = "3.067"
    "E33" = "  +9.41%  "
    "D34" = "4.584"
    "E34" = "  +3.47%  "
    "D35" = "1.035"
    "E35" = "  +2.91%  "
    "D36" = "1.166"
    "E36" = "  +3.81%  "
    "D37" = "0.01989"
    "E37" = "  +3.26%  "
    "D38" = "0.05351"
    "E38" = "  +2.26%  "
    "D39" = "0.5200"
    "E39" = "  +2.03%  "
    "D40" = "2.844"
    "E40" = "  +3.58%  "
    "E41" = "  +2.37%  "
    "D42" = "6.849"
    "E42" = "  +5.64%  "
    "D43" = "8.715"
    "E43" = "  +5.12%  "
    "D44" = "109.61"
    "E44" = "  +1.99%  "
    "E45" = "  +1.76%  "
    "D46" = "1.727"
    "E46" = "  +4.61%  "
    "D47" = "0.4701"
    "E47" = "  +2.55%  "
    "D48" = "0.06429"
    "E48" = "  +2.52%  "
    "D49" = "1.873"
    "E49" = "  +4.06%  "
    "E50" = "  +4.43%  "
    "D51" = "64.55"
    "E51" = "  +1.48%  "
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    # Matches plain decimal numbers (e.g. "15.50", "0.5200", "324.65")
    # but NOT values with more than one '.' (e.g. "27.739.53") and NOT
    # the percent-volume strings (they contain '%', '+', spaces).
    if ($value -match '^-?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.ClearFormats()
    } else {
        $cell.Value = $value
    }
}
